$wb = $excel.ActiveWorkbook

# "2016-09-01 09:23:26" is a shared string used by both Overview!G2 and
# de-de!H2 (same timestamp, coincidentally). Updating it updates both cells.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 09:24:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-01 09:24:19"
$wsDeDe.Range("K2").Value = "2016-09-01 09:24:40"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 09:24:14"
$wsZhCn.Range("K2").Value = "2016-09-01 09:24:33"
